# "expansão das análises automáticas"
# Adds 3 new computed columns (apoio_medio, contribuicoes, media_contribuicoes)
# and rescales the existing "particip"/"taxa_sucesso" columns (E,F) from
# fractions (0-1, 0.00% display) to already-scaled percentage numbers (0-100).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rescale columns E (particip) and F (taxa_sucesso) for rows 2-7 ---
# They keep the same 0.00% number format, but the underlying stored value
# moves from a 0..1 fraction to an already-multiplied-by-100 number.
for ($r = 2; $r -le 7; $r++) {
    $eVal = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 5).Value = $eVal * 100

    $fVal = $ws.Cells.Item($r, 6).Value()
    $ws.Cells.Item($r, 6).Value = $fVal * 100
}

# --- New header cells L1:N1 ---
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Match the header styling used by the rest of row 1 (bold, bordered, centered)
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:N1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- New data columns L, M, N (rows 2-7), default/general number format ---
$apoioMedio = @(91.17116223450478, 91.74154684374953, 90.6205976008235, 88.85439072913162, 18.28712748796549, 28.57643247462115)
$contribuicoes = @(202460, 61093, 151336, 52310, 1935, 273)
$mediaContribuicoes = @(321.8759936406995, 303.9452736318408, 141.7003745318352, 166.0634920634921, 14.33333333333333, 16.05882352941176)

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 12).Value = $apoioMedio[$i]
    $ws.Cells.Item($r, 13).Value = $contribuicoes[$i]
    $ws.Cells.Item($r, 14).Value = $mediaContribuicoes[$i]
}

$wb.Save()
